# "Cambios en feature 1"
#
# The document ends with an empty paragraph whose only content is the
# underline formatting carried on its paragraph mark (w:pPr/w:rPr/w:u).
# We fill that paragraph with the underlined sentence
#   "Funcionalidad agregada en feature 1."
# (with "feature" wrapped the way Word's proofer wraps a flagged word),
# and the paragraph mark's own direct formatting goes away once the
# paragraph actually carries real run-level formatting instead.

$d = $word.ActiveDocument

# The target (empty, underline-only) paragraph is the last one in the body.
$target = $d.Paragraphs.Item($d.Paragraphs.Count)

# Clearing/normalizing the paragraph style up front drops the leftover
# w:pPr/w:rPr (paragraph-mark-only formatting) once the paragraph gets
# real content; the underline on the new text is applied explicitly below.
$target.Style = "Normal"

# Insert the new sentence as proper WordprocessingML so it lands as three
# separate, explicitly-underlined runs with proofing marks bracketing
# "feature" (as Word itself would leave them from its spell-checker),
# instead of one run that happens to render the same text. Anchoring the
# insertion point right after the previous paragraph's own mark (rather
# than at the start of the empty target paragraph) lets the new content
# take over the target paragraph in place instead of pushing it down.
$anchor = $target.Previous()
$insertionPoint = $d.Range($anchor.Range.End, $anchor.Range.End)

$wordOpenXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r>
              <w:rPr>
                <w:u w:val="single"/>
              </w:rPr>
              <w:t xml:space="preserve">Funcionalidad agregada en </w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:rPr>
                <w:u w:val="single"/>
              </w:rPr>
              <w:t>feature</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:rPr>
                <w:u w:val="single"/>
              </w:rPr>
              <w:t xml:space="preserve"> 1.</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$insertionPoint.InsertXML($wordOpenXml)
